$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 743 (shifts the existing rows 743..823 down to 744..824,
# growing the sheet's used range from A1:R823 to A1:R824).
$ws.Rows.Item(743).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Cells.Item(743, 1).Value = 3
$ws.Cells.Item(743, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(743, 3).Value = "Coquimbo"
$ws.Cells.Item(743, 4).Value = 45194
$ws.Cells.Item(743, 5).Value = 5
$ws.Cells.Item(743, 6).Value = 100112037
$ws.Cells.Item(743, 7).Value = "Cebollín"
$ws.Cells.Item(743, 8).Value = "Sin especificar"
$ws.Cells.Item(743, 9).Value = "Primera"
$ws.Cells.Item(743, 10).Value = 120
$ws.Cells.Item(743, 11).Value = 3500
$ws.Cells.Item(743, 12).Value = 3500
$ws.Cells.Item(743, 13).Value = 3500
$ws.Cells.Item(743, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(743, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(743, 16).Value = 97
$ws.Cells.Item(743, 17).Value = 36
$ws.Cells.Item(743, 18).Value = "Hortaliza"
